$wb = $excel.ActiveWorkbook

# The source sheet "Sheet_250222_231013" already has the exact table layout
# (ID / Name / Current Status / Comment, A1:D8) that both new sheets are
# based on, so duplicate it via Copy (keeps formatting/page setup/styles)
# and place the copies at the end of the tab strip.
$src = $wb.Worksheets.Item("Sheet_250222_231013")

# --- New sheet: Sheet_260222_192336 (straight duplicate, A1:D8) ---
$src.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws1 = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws1.Name = "Sheet_260222_192336"

# --- New sheet: Sheet_260222_192617 (duplicate plus an inserted
#     "Comment from Dt1" column before the Comment column, A1:E8) ---
$src.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws2 = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws2.Name = "Sheet_260222_192617"

$ws2.Columns.Item(4).Insert()
$ws2.Range("D1").Value = "Comment from Dt1"
$ws2.Range("D2:D8").Clear()

# Restore the originally active sheet/tab so the rest of the workbook's
# view-state stays untouched.
$wb.Worksheets.Item("Result").Activate()
